$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data (row 4 had a value that needs to disappear), including
# formatting / phonetic info left over from the old content.
$ws.Cells.Clear()

# Set column A first (Item, Impacto1, Impacto2), then column B (peso, 0.3, 0.7)
# to match the shared-string insertion order of the target workbook.
$ws.Range("A1").Value = "Item"
$ws.Range("A2").Value = "Impacto1"
$ws.Range("A3").Value = "Impacto2"

$ws.Range("B1").Value = "peso"
$ws.Range("B2").Value = 0.3
$ws.Range("B3").Value = 0.7

# Set selection to B4 (as seen in the diff)
$ws.Range("B4").Select()
